# Update "想去人数" (want-to-go count) figures in column F on the
# "展览" and "全部类型" sheets. Both sheets carry the same rows of data,
# so the same set of row -> new value updates applies to each.

$wb = $excel.ActiveWorkbook

$updates = @{
    5  = 53
    14 = 1697
    16 = 490
    22 = 1422
    23 = 3354
    27 = 1088
    28 = 80
    29 = 1779
    32 = 55
    38 = 32
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
